$d = $word.ActiveDocument

# The document contains one Word field whose code is the M2Doc expression:
#   m:'Obeo\'s website'.asLink('http://www.obeo.fr', 'This will open the Obeo website.')
# represented with the usual begin/instrText*/end run sequence.
#
# The parser now uses TokenIteratorFieldRewriterSplit, which rewrites such a
# field as plain literal text runs containing the M2Doc template syntax
# (curly braces around the expression) instead of a real Word field, split
# across the same number of runs as the original instrText runs:
#   {m:'Obeo\'s website'.asLink('http://www.obeo.fr'
#   , 'This will open the Obeo website.'
#   )}

$field = $d.Fields.Item(1)

# Locate the paragraph holding the field (without assuming a fixed index) so
# we know where to re-insert the rewritten text, then remove the field
# (begin/instrText/end runs) entirely.
$fieldCodeStart = $field.Code.Start
$insertAt = $fieldCodeStart
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i).Range
    if (($fieldCodeStart -ge $candidate.Start) -and ($fieldCodeStart -le $candidate.End)) {
        $insertAt = $candidate.Start
        break
    }
}
$field.Delete()

$target = $d.Range($insertAt, $insertAt)

# Build the replacement runs as literal OOXML so each chunk of text keeps its
# own <w:r>, mirroring the 3 runs produced by the rewriter (instead of being
# auto-coalesced into a single run because they all share the same rPr).
$run1 = "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>{m:'Obeo\'s website'.asLink('http://www.obeo.fr'</w:t></w:r>"
$run2 = "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>, 'This will open the Obeo website.'</w:t></w:r>"
$run3 = "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">)}</w:t></w:r>"

$body = "<w:p>$run1$run2$run3</w:p>"

$package = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
  "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
  "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
  "<w:body>$body</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$target.InsertXML($package)
